$d = $word.ActiveDocument

function Replace-InParagraph($Index, $OldText, $NewText) {
    $p = $d.Paragraphs.Item($Index)
    $rng = $p.Range.Duplicate
    $ok = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $ok) {
        throw "Replace failed in paragraph $Index for text: $OldText"
    }
}

function Find-LabelRange($Label, $SearchStart, $SearchEnd) {
    $rng = $d.Range($SearchStart, $SearchEnd)
    $ok = $rng.Find.Execute($Label, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) { throw "label not found: $Label" }
    $result = @($rng.Start, $rng.End)
    return $result
}

function Get-ParaBounds($Index) {
    $p = $d.Paragraphs.Item($Index)
    $result = @($p.Range.Start, $p.Range.End)
    return $result
}

# --- Paragraph 6: "Objetivos" (PT) body ---
Replace-InParagraph 6 "Propiciar aos alunos conhecimento sobre os fundamentos, objetivos e métodos da Avaliação de Impacto Ambiental." "Os objetivos da Avaliação de Impacto Ambiental (AIA). O processo da AIA: triagem, definição do escopo, estudos de base, análise de impactos ambientais, mitigação, análise técnica e acompanhamento. Requisitos legais."

# --- Paragraph 7: "Objetivos" (EN, italic) body ---
Replace-InParagraph 7 "Provide knowledge on the fundamentals, objectives and methods of environmental impact assessment." "Environmental impact assessment process and its objectives; Environmental impacts on surface water, groundwater and ocean; Water pollution; Environmental impacts on soil; Atmospheric pollution and human health."

# --- Paragraph 9: Docente(s) list item ---
Replace-InParagraph 9 "5464150 - Mariana Consiglio Kasemodel" "Propiciar aos alunos conhecimento sobre os fundamentos, objetivos e métodos da Avaliação de Impacto Ambiental."

# --- Paragraph 11: "Programa resumido" (PT) body ---
Replace-InParagraph 11 "Os objetivos da Avaliação de Impacto Ambiental (AIA). O processo da AIA: triagem, definição do escopo, estudos de base, análise de impactos ambientais, mitigação, análise técnica e acompanhamento. Requisitos legais." "Conceitos básicos e definições. Origem e difusão da AIA. AIA e licenciamento: objetivos e fundamentos. Quadro legal e institucional brasileiro para a AIA. O processo de AIA e seus componentes. Etapas do planejamento e execução de um Estudo de Impacto Ambiental. Alternativas tecnológicas e de localização. Estudos de base e diagnóstico ambiental. Técnicas de identificação e previsão de impactos. Métodos e critérios de avaliação da importância dos impactos. Plano de gestão ambiental: medidas mitigadoras, compensatórias, de valorização e monitoramento. Tomada de decisão e acompanhamento. Estudos de caso. A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina."

# --- Paragraph 12: "Programa resumido" (EN, italic) body ---
Replace-InParagraph 12 "Environmental impact assessment process and its objectives; Environmental impacts on surface water, groundwater and ocean; Water pollution; Environmental impacts on soil; Atmospheric pollution and human health." "Provide knowledge on the fundamentals, objectives and methods of environmental impact assessment."

# --- Paragraph 14: "Programa" (EN, italic) body ---
Replace-InParagraph 14 "Conceitos básicos e definições. Origem e difusão da AIA. AIA e licenciamento: objetivos e fundamentos. Quadro legal e institucional brasileiro para a AIA. O processo de AIA e seus componentes. Etapas do planejamento e execução de um Estudo de Impacto Ambiental. Alternativas tecnológicas e de localização. Estudos de base e diagnóstico ambiental. Técnicas de identificação e previsão de impactos. Métodos e critérios de avaliação da importância dos impactos. Plano de gestão ambiental: medidas mitigadoras, compensatórias, de valorização e monitoramento. Tomada de decisão e acompanhamento. Estudos de caso. A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina." "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."

# --- Paragraph 17: "Avaliação" (Método / Critério / Norma de recuperação list) ---
# Shift the three value runs one slot earlier, then append the bibliography block
# after "Norma de recuperação: ".

# STEP 1: replace value after "Método: "
$b = Get-ParaBounds 17
$metodo = Find-LabelRange "Método: " $b[0] $b[1]
$criterio = Find-LabelRange "Critério: " $metodo[1] $b[1]
$value1 = $d.Range($metodo[1], $criterio[0])
$ok1 = $value1.Find.Execute("Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas.", $true, $false, $false, $false, $false, $true, 1, $false, "Média ponderada de atividades e provas.", 2)
if (-not $ok1) { throw "P17 step1 replace failed" }

# STEP 2: replace value after "Critério: "
$b = Get-ParaBounds 17
$metodo = Find-LabelRange "Método: " $b[0] $b[1]
$criterio = Find-LabelRange "Critério: " $metodo[1] $b[1]
$norma = Find-LabelRange "Norma de recuperação: " $criterio[1] $b[1]
$value2 = $d.Range($criterio[1], $norma[0])
$ok2 = $value2.Find.Execute("Média ponderada de atividades e provas.", $true, $false, $false, $false, $false, $true, 1, $false, "1 (uma) prova escrita", 2)
if (-not $ok2) { throw "P17 step2 replace failed" }

# STEP 3: replace value after "Norma de recuperação: " with the bibliography block
$b = Get-ParaBounds 17
$metodo = Find-LabelRange "Método: " $b[0] $b[1]
$criterio = Find-LabelRange "Critério: " $metodo[1] $b[1]
$norma = Find-LabelRange "Norma de recuperação: " $criterio[1] $b[1]
$value3 = $d.Range($norma[1], $b[1])
$newBib = "Bibliografia básica^lSÁNCHEZ, L.E., Avaliação de impacto ambiental: conceitos e métodos. Oficina de textos: São Paulo, 2013. 583p.^lCALIJURI, M.C., CUNHA, D.G.F. (Org.), Engenharia ambiental: conceitos, tecnologia e gestão. Elsevier: Rio de Janeiro, 2019. 685p.^l^lBibliografia complementar:^lCOMPANHIA AMBIENTAL DO ESTADO DE SÃO PAULO - CETESB. MANUAL PARA ELABORAÇÃO DE ESTUDOS PARA O LICENCIAMENTO COM AVALIAÇÃO DE IMPACTO AMBIENTAL. Departamento de Desenvolvimento de Ações Estratégicas para o Licenciamento da Diretoria I- ID - CETESB. Anexo único, 2014. 250p."
$ok3 = $value3.Find.Execute("1 (uma) prova escrita", $true, $false, $false, $false, $false, $true, 1, $false, $newBib, 2)
if (-not $ok3) { throw "P17 step3 replace failed" }

# --- Paragraph 19: Bibliografia body -> replaced by docente list item ---
# This paragraph's whole content (several runs/breaks) collapses into a single
# plain run, so replace the entire paragraph range (excluding the paragraph
# mark) rather than doing a partial text Find/Replace.
$p19 = $d.Paragraphs.Item(19)
$rng19 = $p19.Range.Duplicate
[void]$rng19.MoveEnd(1, -1)
$rng19.Text = "5464150 - Mariana Consiglio Kasemodel"
